# Cash Flow Quantization Size.xlsx - bring CFQS sheet up to date:
#  - label in B1 gains a "($)" units suffix
#  - quantization size value in B2 doubles (50000 -> 100000)
#  - column B widened slightly to fit the new label text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CFQS")

$ws.Range("B1").Value = "Quantization Size ($)"
$ws.Range("B2").Value = 100000
$ws.Columns.Item(2).ColumnWidth = 18.666666666666668
